$d = $word.ActiveDocument

# --- Step 1: remove the "License Information" heading paragraph (directly follows the
#     "Resource: ..." Heading 2 paragraph). Range.Delete() on a paragraph's full range
#     removes the paragraph (including its mark), merging it away. ---
$pLic = $d.Paragraphs.Item(4)
$pLic.Range.Delete()

# --- Step 2: clear the attribution paragraph's content (now paragraph 4) but keep the
#     paragraph itself - trim the trailing paragraph-mark unit off the range first so the
#     delete does not merge this paragraph into the next one. This also removes the two
#     hyperlink runs (Biblica, Inc. / CC BY-SA 4.0 license) that lived in this paragraph. ---
$p4 = $d.Paragraphs.Item(4)
$r4 = $p4.Range
$r4.End = $r4.End - 1
$r4.Delete()

# --- Step 3: remove the "This PDF version is provided under the same license." paragraph
#     (now paragraph 5) entirely. ---
$pPdf = $d.Paragraphs.Item(5)
$pPdf.Range.Delete()

# --- Step 4: insert the new attribution text into the now-empty paragraph 4. ---
$p4b = $d.Paragraphs.Item(4)
$r4b = $p4b.Range
$r4b.InsertAfter("Biblica Study Notes (Key Terms) © 2023 Biblica Inc. Released under CC BY-SA 4.0 license. Biblica Study Notes has been adapted in the following languages: Tok Pisin, Arabic (عربي), French (Français), Hindi (हिंदी), Indonesian (Bahasa Indonesia), Portuguese (Português), Russian (Русский), Spanish (Español), Swahili (Kiswahili), and Simplified Chinese (简体中文)from Biblica Study Notes © 2023 Biblica Inc. Released under CC BY-SA 4.0 license by Mission Mutual.")

# --- Step 5: bold the "Biblica Study Notes (Key Terms)" lead-in only. ---
$p4c = $d.Paragraphs.Item(4)
$r4c = $p4c.Range
$r4c.End = $r4c.Start + 31
$r4c.Bold = 1

# --- Step 6: remove the italic "牧羊人" paragraph that follows the "mu"
#     Heading 2 paragraph (a Normal-style paragraph; there is a later Heading-2
#     paragraph with the same text that must be left untouched). ---
foreach ($p in $d.Paragraphs) {
    $txt = $p.Range.Text.TrimEnd([char]13)
    if (($txt -eq "牧羊人") -and ($p.Style.NameLocal -ne "Heading 2")) {
        $p.Range.Delete()
        break
    }
}

Write-Output "DONE"
